# Applies the periodic "cryptos" price/volume refresh produced by the
# GitHub Actions scraper job. Cell values are written as literal text
# (matching the inlineStr cells already in the sheet); a leading
# apostrophe forces Excel to keep purely-numeric-looking strings (e.g.
# "9.60", "0.0760") as text instead of silently recasting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '44.185.06'
$ws.Range('E2').Value = '  +2.46%  '
# Row 3
$ws.Range('D3').Value = '2.425.98'
$ws.Range('E3').Value = '  +2.04%  '
# Row 4
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.08%  '
# Row 5
$ws.Range('D5').Value = '''307.95'
$ws.Range('E5').Value = '  +1.64%  '
# Row 6
$ws.Range('D6').Value = '''100.97'
$ws.Range('E6').Value = '  +4.12%  '
# Row 7
$ws.Range('D7').Value = '''0.512'
$ws.Range('E7').Value = '  +1.39%  '
# Row 8
$ws.Range('E8').Value = '  -0.05%  '
# Row 9
$ws.Range('E9').Value = '  +0.42%  '
# Row 10
$ws.Range('D10').Value = '''35.23'
$ws.Range('E10').Value = '  +2.95%  '
# Row 11
$ws.Range('E11').Value = '  +1.73%  '
# Row 12
$ws.Range('D12').Value = '''19.04'
$ws.Range('E12').Value = '  +3.91%  '
# Row 13
$ws.Range('E13').Value = '  +2.03%  '
# Row 14
$ws.Range('E14').Value = '  +1.94%  '
# Row 15
$ws.Range('D15').Value = '2.804.47'
$ws.Range('E15').Value = '  +1.94%  '
# Row 16
$ws.Range('D16').Value = '2.430.06'
$ws.Range('E16').Value = '  +3.36%  '
# Row 17
$ws.Range('D17').Value = '''0.836'
$ws.Range('E17').Value = '  +3.67%  '
# Row 18
$ws.Range('D18').Value = '44.122.79'
$ws.Range('E18').Value = '  +2.26%  '
# Row 19
$ws.Range('E19').Value = '  +1.08%  '
# Row 20
$ws.Range('D20').Value = '''6.41'
$ws.Range('E20').Value = '  +1.97%  '
# Row 21
$ws.Range('E21').Value = '  +1.95%  '
# Row 22
$ws.Range('D22').Value = '''68.58'
$ws.Range('E22').Value = '  +0.24%  '
# Row 23
$ws.Range('B23').Value = 'ImmutableX'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D23').Value = '''2.32'
$ws.Range('E23').Value = '  +5.12%  '
# Row 24
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '''240.42'
$ws.Range('E24').Value = '  +2.21%  '
# Row 25
$ws.Range('E25').Value = '  +1.64%  '
# Row 26
$ws.Range('E26').Value = '  -0.05%  '
# Row 27
$ws.Range('D27').Value = '''25.26'
$ws.Range('E27').Value = '  +1.65%  '
# Row 28
$ws.Range('E28').Value = '  -1.56%  '
# Row 29
$ws.Range('D29').Value = '''9.60'
$ws.Range('E29').Value = '  +5.07%  '
# Row 30
$ws.Range('D30').Value = '''32.86'
$ws.Range('E30').Value = '  +4.82%  '
# Row 31
$ws.Range('D31').Value = '''18.69'
$ws.Range('E31').Value = '  +7.47%  '
# Row 32
$ws.Range('E32').Value = '  +11.15%  '
# Row 33
$ws.Range('D33').Value = '''5.19'
$ws.Range('E33').Value = '  +2.32%  '
# Row 34
$ws.Range('E34').Value = '  -0.10%  '
# Row 35
$ws.Range('D35').Value = '''0.0760'
$ws.Range('E35').Value = '  +2.01%  '
# Row 36
$ws.Range('D36').Value = '''1.89'
$ws.Range('E36').Value = '  +3.20%  '
# Row 37
$ws.Range('D37').Value = '''4.48'
$ws.Range('E37').Value = '  +4.29%  '
# Row 38
$ws.Range('D38').Value = '''129.73'
$ws.Range('E38').Value = '  +25.54%  '
# Row 39
$ws.Range('E39').Value = '  +3.70%  '
# Row 40
$ws.Range('E40').Value = '  -0.84%  '
# Row 41
$ws.Range('E41').Value = '  +0.80%  '
# Row 42
$ws.Range('D42').Value = '''21.36'
$ws.Range('E42').Value = '  -5.12%  '
# Row 43
$ws.Range('E43').Value = '  +2.50%  '
# Row 44
$ws.Range('D44').Value = '1.951.32'
$ws.Range('E44').Value = '  -0.53%  '
# Row 45
$ws.Range('E45').Value = '  +1.87%  '
# Row 46
$ws.Range('E46').Value = '  +4.88%  '
# Row 47
$ws.Range('D47').Value = '''9.43'
$ws.Range('E47').Value = '  +3.40%  '
# Row 48
$ws.Range('D48').Value = '''1.64'
$ws.Range('E48').Value = '  +9.19%  '
# Row 49
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.660.93'
$ws.Range('E49').Value = '  +2.21%  '
# Row 50
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '''53.49'
$ws.Range('E50').Value = '  +1.67%  '
# Row 51
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '''73.80'
$ws.Range('E51').Value = '  +2.56%  '
